$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right before current row 18 (shifts old rows 18..106 down to 19..107)
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly entry
$ws.Cells.Item(18,1).Value2  = 5
$ws.Cells.Item(18,2).Value   = "Macroferia Regional de Talca"
$ws.Cells.Item(18,3).Value   = "Maule"
$ws.Cells.Item(18,4).Value2  = 44831
$ws.Cells.Item(18,5).Value2  = 7
$ws.Cells.Item(18,6).Value2  = 100112013
$ws.Cells.Item(18,7).Value   = "Alcachofa"
$ws.Cells.Item(18,8).Value   = "Madrigal"
$ws.Cells.Item(18,9).Value   = "Primera"
$ws.Cells.Item(18,10).Value2 = 400
$ws.Cells.Item(18,11).Value2 = 10000
$ws.Cells.Item(18,12).Value2 = 10000
$ws.Cells.Item(18,13).Value2 = 10000
$ws.Cells.Item(18,14).Value  = "`$/caja 40 unidades"
$ws.Cells.Item(18,15).Value  = "Provincia del Elquí"
$ws.Cells.Item(18,16).Value2 = 250
$ws.Cells.Item(18,17).Value2 = 40
$ws.Cells.Item(18,18).Value  = "Hortaliza"

# Append a brand-new row 108 with another new weekly entry
$ws.Cells.Item(108,1).Value2  = 5
$ws.Cells.Item(108,2).Value   = "Macroferia Regional de Talca"
$ws.Cells.Item(108,3).Value   = "Maule"
$ws.Cells.Item(108,4).Value2  = 44832
$ws.Cells.Item(108,5).Value2  = 7
$ws.Cells.Item(108,6).Value2  = 100112013
$ws.Cells.Item(108,7).Value   = "Alcachofa"
$ws.Cells.Item(108,8).Value   = "Madrigal"
$ws.Cells.Item(108,9).Value   = "Primera"
$ws.Cells.Item(108,10).Value2 = 300
$ws.Cells.Item(108,11).Value2 = 10000
$ws.Cells.Item(108,12).Value2 = 10000
$ws.Cells.Item(108,13).Value2 = 10000
$ws.Cells.Item(108,14).Value  = "`$/caja 40 unidades"
$ws.Cells.Item(108,15).Value  = "Provincia del Elquí"
$ws.Cells.Item(108,16).Value2 = 250
$ws.Cells.Item(108,17).Value2 = 40
$ws.Cells.Item(108,18).Value  = "Hortaliza"

# Make sure the date column keeps its original format style ("D" column) for both new rows
$ws.Cells.Item(18,4).NumberFormat  = $ws.Cells.Item(19,4).NumberFormat
$ws.Cells.Item(108,4).NumberFormat = $ws.Cells.Item(19,4).NumberFormat

# Update the worksheet dimension / used range reflects new data automatically
